$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E3 value (-5.7) is removed -> becomes blank/missing
$ws.Range("E3").Value = ""

# Row 26 ("RM 232") is removed entirely from the dataset; all subsequent
# rows shift up by one.
$ws.Rows.Item(26).Delete()

# After the shift, the row formerly holding "SC 92" is now row 27; it is
# also removed entirely, shifting everything below it up by one more.
$ws.Rows.Item(27).Delete()

# Now rows 26-33 hold (in order): SC 5, SC 101, SC 105, SC 119, SC 120,
# SC 132, SC 193, SC 232. Update the individual cells that changed value.

# SC 5 (row 26): column D was missing, now imputed to -13.8
$ws.Range("D26").Value = -13.8

# SC 101 (row 27): column D was -14.6, now removed/missing
$ws.Range("D27").Value = ""

# SC 232 (row 33): columns D and E were missing, now imputed
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
